$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (Mac-Addresses) following the existing pattern
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 10030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 10031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

# Update the view: scroll back to top-left and select E31
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
